# FDR RQ3 code update
# Re-ran FDR corrections for RQ3 (p-value column G) on all four sheets.
# Updated p-values, re-applied the corrected-pvalue number font/style, and
# cleared the "*" significance marker (shared string 13) on rows that are
# no longer significant after the re-run. Also refreshes sheet selections.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # CCL5
$ws2 = $wb.Worksheets.Item(2)   # sqrt(CCL5)
$ws3 = $wb.Worksheets.Item(3)   # CLSTN3
$ws4 = $wb.Worksheets.Item(4)   # NEGR1

# ---------------------------------------------------------------------
# 1. Updated FDR-corrected p-values (column G, rows 2-10) per sheet
# ---------------------------------------------------------------------

$ws1.Range("G2").Value = 0.86368551089373902
$ws1.Range("G3").Value = 0.39259232958547502
$ws1.Range("G4").Value = 0.64850294779134998
$ws1.Range("G5").Value = 0.64850294779134998
$ws1.Range("G6").Value = 0.64850294779134998
$ws1.Range("G7").Value = 0.28925150014732798
$ws1.Range("G8").Value = 0.89473323361632695
$ws1.Range("G9").Value = 0.89473323361632695
$ws1.Range("G10").Value = 0.64850294779134998

$ws2.Range("G2").Value = 0.98860546879340805
$ws2.Range("G3").Value = 0.0994878509188923
$ws2.Range("G4").Value = 0.75471628840566696
$ws2.Range("G5").Value = 0.39848108508763502
$ws2.Range("G6").Value = 0.98860546879340805
$ws2.Range("G7").Value = 0.39848108508763502
$ws2.Range("G8").Value = 0.75471628840566696
$ws2.Range("G9").Value = 0.98860546879340805
$ws2.Range("G10").Value = 0.103564360938349

$ws3.Range("G2").Value = 0.36163325946404101
$ws3.Range("G3").Value = 0.88083698415301703
$ws3.Range("G4").Value = 0.16465098547288801
$ws3.Range("G5").Value = 0.16465098547288801
$ws3.Range("G6").Value = 0.16465098547288801
$ws3.Range("G7").Value = 0.16465098547288801
$ws3.Range("G8").Value = 0.50334521700138202
$ws3.Range("G9").Value = 0.21259401715891399
$ws3.Range("G10").Value = 0.667342603653165

$ws4.Range("G2").Value = 0.28203286864166
$ws4.Range("G3").Value = 0.83910256252064996
$ws4.Range("G4").Value = 0.307102741296917
$ws4.Range("G5").Value = 0.37292449368537001
$ws4.Range("G6").Value = 0.48305832476855898
$ws4.Range("G7").Value = 0.307102741296917
$ws4.Range("G8").Value = 0.83910256252064996
$ws4.Range("G9").Value = 0.48305832476855898
$ws4.Range("G10").Value = 0.307102741296917

# ---------------------------------------------------------------------
# 2. Re-style the recalculated p-value column with an explicit black font
#    (was inheriting the theme text color; now pinned to RGB black) on
#    all four sheets.
# ---------------------------------------------------------------------

$ws1.Range("G2:G10").Font.Color = 0
$ws2.Range("G2:G10").Font.Color = 0
$ws3.Range("G2:G10").Font.Color = 0
$ws4.Range("G2:G10").Font.Color = 0

# ---------------------------------------------------------------------
# 3. Rows that are no longer significant after the re-run lose their "*"
#    marker in column H (CLSTN3 row 5, NEGR1 row 2). The other two
#    significance markers (CCL5 row 7, sqrt(CCL5) rows 3 & 10) remain.
# ---------------------------------------------------------------------

$ws3.Range("H5").ClearContents()
$ws4.Range("H2").ClearContents()

# ---------------------------------------------------------------------
# 4. Refresh sheet selections / active sheet to reflect where the author
#    was working after the re-run. Sheet1 is activated last so it ends
#    up as the active/selected tab.
# ---------------------------------------------------------------------

$ws2.Range("G2:G10").Select()
$ws3.Range("G2:G10").Select()
$ws4.Range("E10").Select()
$ws1.Range("F12").Select()
